$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing rows 2-18: new WIID (column A) and Date (column E) values
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "95310357"
$ws.Range("A2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2022-05-02"
$ws.Range("E2").Style = "Normal"
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "95310351"
$ws.Range("A3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2022-06-04"
$ws.Range("E3").Style = "Normal"
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "95310347"
$ws.Range("A4").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2021-04-07"
$ws.Range("E4").Style = "Normal"
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "95310364"
$ws.Range("A5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2020-11-22"
$ws.Range("E5").Style = "Normal"
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "95310362"
$ws.Range("A6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2020-07-21"
$ws.Range("E6").Style = "Normal"
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "95310358"
$ws.Range("A7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2019-09-07"
$ws.Range("E7").Style = "Normal"
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "95310356"
$ws.Range("A8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2018-05-20"
$ws.Range("E8").Style = "Normal"
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "95310348"
$ws.Range("A9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2018-08-12"
$ws.Range("E9").Style = "Normal"
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "95310355"
$ws.Range("A10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2017-09-23"
$ws.Range("E10").Style = "Normal"
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "95310363"
$ws.Range("A11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2023-07-29"
$ws.Range("E11").Style = "Normal"
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "95310349"
$ws.Range("A12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2019-01-08"
$ws.Range("E12").Style = "Normal"
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "95310352"
$ws.Range("A13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2023-03-20"
$ws.Range("E13").Style = "Normal"
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "95310360"
$ws.Range("A14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2020-02-23"
$ws.Range("E14").Style = "Normal"
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "95310346"
$ws.Range("A15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2019-04-08"
$ws.Range("E15").Style = "Normal"
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "95310359"
$ws.Range("A16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2023-01-04"
$ws.Range("E16").Style = "Normal"
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = "95310361"
$ws.Range("A17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2021-04-16"
$ws.Range("E17").Style = "Normal"
$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = "95310350"
$ws.Range("A18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2017-12-21"
$ws.Range("E18").Style = "Normal"

# Add new rows 19 and 20
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = "95310354"
$ws.Range("A19").Style = "Normal"
$ws.Range("B19").Value = "Calculate Client Security Hash"
$ws.Range("C19").Value = "WI5"
$ws.Range("D19").Value = "Open"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2019-11-28"
$ws.Range("E19").Style = "Normal"
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = "95310353"
$ws.Range("A20").Style = "Normal"
$ws.Range("B20").Value = "Calculate Client Security Hash"
$ws.Range("C20").Value = "WI5"
$ws.Range("D20").Value = "Open"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2021-08-25"
$ws.Range("E20").Style = "Normal"
